$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "tons_per_lanemi" column (C) ---

# Copy the number formatting from column B (Tons) down onto column C for the
# existing data rows (2-61) so the new cells share the same number style.
$ws.Range("B2:B61").Copy()
$ws.Range("C2:C61").PasteSpecial(-4122)

# Data values for the new column (tons of road salt per lane-mile)
$ws.Cells.Item(2, 3).Value = 4.8
$ws.Cells.Item(3, 3).Value = 2.7
$ws.Cells.Item(4, 3).Value = 5.5
$ws.Cells.Item(5, 3).Value = 3.9
$ws.Cells.Item(6, 3).Value = 4.2
$ws.Cells.Item(7, 3).Value = 7.5
$ws.Cells.Item(8, 3).Value = 5.7
$ws.Cells.Item(9, 3).Value = 8
$ws.Cells.Item(10, 3).Value = 6.2
$ws.Cells.Item(11, 3).Value = 8.5
$ws.Cells.Item(12, 3).Value = 8.7
$ws.Cells.Item(13, 3).Value = 11.8
$ws.Cells.Item(14, 3).Value = 8.7
$ws.Cells.Item(15, 3).Value = 10
$ws.Cells.Item(16, 3).Value = "N/A"
$ws.Cells.Item(17, 3).Value = "N/A"
$ws.Cells.Item(18, 3).Value = "N/A"
$ws.Cells.Item(19, 3).Value = "N/A"
$ws.Cells.Item(20, 3).Value = "N/A"
$ws.Cells.Item(21, 3).Value = "N/A"
$ws.Cells.Item(22, 3).Value = "N/A"
$ws.Cells.Item(23, 3).Value = "N/A"
$ws.Cells.Item(24, 3).Value = "N/A"
$ws.Cells.Item(25, 3).Value = 8.6
$ws.Cells.Item(26, 3).Value = 8.2
$ws.Cells.Item(27, 3).Value = 7.9
$ws.Cells.Item(28, 3).Value = 11
$ws.Cells.Item(29, 3).Value = 7.1
$ws.Cells.Item(30, 3).Value = 8.1
$ws.Cells.Item(31, 3).Value = 8.3
$ws.Cells.Item(32, 3).Value = 10.6
$ws.Cells.Item(33, 3).Value = 13
$ws.Cells.Item(34, 3).Value = 12
$ws.Cells.Item(35, 3).Value = 14.8
$ws.Cells.Item(36, 3).Value = 11.1
$ws.Cells.Item(37, 3).Value = 10.4
$ws.Cells.Item(38, 3).Value = 15.5
$ws.Cells.Item(39, 3).Value = 17.8
$ws.Cells.Item(40, 3).Value = 14
$ws.Cells.Item(41, 3).Value = 12.4
$ws.Cells.Item(42, 3).Value = 11.4
$ws.Cells.Item(43, 3).Value = 17.1
$ws.Cells.Item(44, 3).Value = 10
$ws.Cells.Item(45, 3).Value = 10.6
$ws.Cells.Item(46, 3).Value = 12.4
$ws.Cells.Item(47, 3).Value = 12.8
$ws.Cells.Item(48, 3).Value = 12.4
$ws.Cells.Item(49, 3).Value = 12.2
$ws.Cells.Item(50, 3).Value = 19.4
$ws.Cells.Item(51, 3).Value = 17
$ws.Cells.Item(52, 3).Value = 12.2
$ws.Cells.Item(53, 3).Value = 17
$ws.Cells.Item(54, 3).Value = 10.5
$ws.Cells.Item(55, 3).Value = 18.2
$ws.Cells.Item(56, 3).Value = 19.5
$ws.Cells.Item(57, 3).Value = 11.3
$ws.Cells.Item(58, 3).Value = 11.6
$ws.Cells.Item(59, 3).Value = 15.2
$ws.Cells.Item(60, 3).Value = 16.4
$ws.Cells.Item(61, 3).Value = 15.9

# Header for the new column (added after the "N/A" cells so the shared
# string table ends up ordered the same way as the source workbook)
$ws.Range("C1").Value = "tons_per_lanemi"

# --- Add a new row for year 2019 (row 62), matching the formatting of the
#     previous row (2018, row 61) ---
$ws.Range("A61:C61").Copy()
$ws.Range("A62:C62").PasteSpecial(-4122)
$ws.Cells.Item(62, 1).Value = 2019
# Column A has no explicit "General" style of its own here (the plain
# number style used by the rest of the year column collapses back to the
# sheet/column default when round-tripped), and the column's own default
# style is a long-date format -- so without an explicit override the new
# year would render as a date instead of the number 2019. Force a plain
# numeric format to keep it displaying as a year, same as the rows above.
$ws.Cells.Item(62, 1).NumberFormat = "0"
$ws.Cells.Item(62, 2).Value = 425558
$ws.Cells.Item(62, 3).Value = 12.2

# --- Match the saved view/selection state ---
$ws.Range("G59").Select()
